$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheets: rename the existing sheet and add the second one
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vue結構"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "問題解決"

# ------------------------------------------------------------------
# Sheet "vue結構" content
# ------------------------------------------------------------------

# Header row - bold, centered
$c = $ws1.Range("B1")
$c.Value = "View"
$c.Font.Name = "微軟正黑體"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

# Numbered list rows (A: index number, B: label - centered, normal weight).
# Note: row 8's label ("cartList") is written to the shared-string table
# before row 7's ("detail") so the shared-string indices line up with the
# target workbook (cartList=6, detail=7) while the row order on the sheet
# stays 1..7 in column A / rows 2..8.
$rows = @(
    @(2, 1, "home"),
    @(3, 2, "men"),
    @(4, 3, "women"),
    @(5, 4, "kids"),
    @(6, 5, "loggin/regist"),
    @(8, 7, "cartList"),
    @(7, 6, "detail")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $idx = $r[1]
    $label = $r[2]

    $a = $ws1.Range("A$rowNum")
    $a.Value = $idx
    $a.Font.Name = "微軟正黑體"

    $b = $ws1.Range("B$rowNum")
    $b.Value = $label
    $b.Font.Name = "微軟正黑體"
    $b.HorizontalAlignment = -4108
    $b.VerticalAlignment = -4108
}

# Second header / footer label - bold, centered
$c2 = $ws1.Range("B11")
$c2.Value = "Nav、Footer"
$c2.Font.Name = "微軟正黑體"
$c2.Font.Bold = $true
$c2.HorizontalAlignment = -4108
$c2.VerticalAlignment = -4108

# Row heights for the populated rows
$ws1.Range("A1:B8").RowHeight = 16.5
$ws1.Range("A11:B11").RowHeight = 16.5

# Column widths
$ws1.Columns.Item(2).ColumnWidth = 18.140625
$ws1.Columns.Item(3).ColumnWidth = 22
$ws1.Columns.Item(4).ColumnWidth = 19.7109375
$ws1.Columns.Item(5).ColumnWidth = 28.85546875

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

[void]$ws1.Range("B15").Select()

# ------------------------------------------------------------------
# Sheet "問題解決" content
# ------------------------------------------------------------------
$d = $ws2.Range("A1")
$d.Value = "子選項的動畫效果：透過設置 max-height 的 transition "
$d.Font.Name = "微軟正黑體"

$ws2.Columns.Item(1).ColumnWidth = 17.5703125

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

[void]$ws2.Range("D9").Select()
